# Auto-generated edit script: update '想去人数' (F column) counts
# across all four worksheets per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 367  # F4: 366 -> 367
$ws.Cells.Item(6, 6).Value = 799  # F6: 796 -> 799
$ws.Cells.Item(9, 6).Value = 2620  # F9: 2617 -> 2620
$ws.Cells.Item(11, 6).Value = 659  # F11: 574 -> 659
$ws.Cells.Item(13, 6).Value = 2043  # F13: 2042 -> 2043
$ws.Cells.Item(15, 6).Value = 6485  # F15: 6481 -> 6485
$ws.Cells.Item(17, 6).Value = 1229  # F17: 1228 -> 1229
$ws.Cells.Item(19, 6).Value = 1481  # F19: 1477 -> 1481
$ws.Cells.Item(21, 6).Value = 1197  # F21: 1194 -> 1197
$ws.Cells.Item(23, 6).Value = 2307  # F23: 2296 -> 2307
$ws.Cells.Item(25, 6).Value = 734  # F25: 731 -> 734
$ws.Cells.Item(27, 6).Value = 5304  # F27: 5303 -> 5304
$ws.Cells.Item(31, 6).Value = 3712  # F31: 3711 -> 3712
$ws.Cells.Item(35, 6).Value = 160  # F35: 159 -> 160

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(22, 6).Value = 243  # F22: 242 -> 243
$ws.Cells.Item(30, 6).Value = 302  # F30: 301 -> 302
$ws.Cells.Item(31, 6).Value = 40  # F31: 39 -> 40
$ws.Cells.Item(32, 6).Value = 134  # F32: 133 -> 134
$ws.Cells.Item(38, 6).Value = 198  # F38: 197 -> 198

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(10, 6).Value = 2785  # F10: 2786 -> 2785
$ws.Cells.Item(11, 6).Value = 289  # F11: 288 -> 289
$ws.Cells.Item(13, 6).Value = 571  # F13: 570 -> 571
$ws.Cells.Item(14, 6).Value = 1159  # F14: 1160 -> 1159

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 367  # F4: 366 -> 367
$ws.Cells.Item(6, 6).Value = 2785  # F6: 2786 -> 2785
$ws.Cells.Item(8, 6).Value = 799  # F8: 796 -> 799
$ws.Cells.Item(11, 6).Value = 2620  # F11: 2617 -> 2620
$ws.Cells.Item(13, 6).Value = 659  # F13: 574 -> 659
$ws.Cells.Item(15, 6).Value = 2043  # F15: 2042 -> 2043
$ws.Cells.Item(17, 6).Value = 6485  # F17: 6481 -> 6485
$ws.Cells.Item(20, 6).Value = 1229  # F20: 1228 -> 1229
$ws.Cells.Item(21, 6).Value = 571  # F21: 570 -> 571
$ws.Cells.Item(22, 6).Value = 1481  # F22: 1477 -> 1481
$ws.Cells.Item(24, 6).Value = 1197  # F24: 1194 -> 1197
$ws.Cells.Item(25, 6).Value = 2307  # F25: 2296 -> 2307
$ws.Cells.Item(29, 6).Value = 734  # F29: 731 -> 734
$ws.Cells.Item(31, 6).Value = 5305  # F31: 5303 -> 5305
$ws.Cells.Item(34, 6).Value = 3712  # F34: 3711 -> 3712
$ws.Cells.Item(35, 6).Value = 302  # F35: 301 -> 302
$ws.Cells.Item(38, 6).Value = 160  # F38: 159 -> 160
$ws.Cells.Item(48, 6).Value = 198  # F48: 197 -> 198
$ws.Cells.Item(49, 6).Value = 198  # F49: 197 -> 198
